$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cluster label ("ECs") added alongside the existing "M2" cluster, following
# Dr Hou's advice. Row 1 headers are unchanged; rows 2-5 now enumerate every
# Sending-cluster x Target-cluster combination for the Cd86/Cd28 ligand-receptor
# pair across the two clusters (ECs, M2).

# Row 2: ECs -> Cd86/Cd28 -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cd86"
$ws.Range("C2").Value = "Cd28"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 8.516883333333334
$ws.Range("H2").Value = 25.55065
$ws.Range("I2").Value = 0.08473352333057485
$ws.Range("J2").Value = 0.08473352333057484
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1791996666666667
$ws.Range("N2").Value = 0.537599
$ws.Range("O2").Value = 0.04251079199666429
$ws.Range("P2").Value = 0.04251079199666429
$ws.Range("Q2").Value = 1.526222654372222
$ws.Range("R2").Value = 13.73600388935
$ws.Range("S2").Value = 0.003602089185450568
$ws.Range("T2").Value = 0.003602089185450568

# Row 3: ECs -> Cd86/Cd28 -> M2
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cd86"
$ws.Range("C3").Value = "Cd28"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 8.516883333333334
$ws.Range("H3").Value = 25.55065
$ws.Range("I3").Value = 0.08473352333057485
$ws.Range("J3").Value = 0.08473352333057484
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.036192666666667
$ws.Range("N3").Value = 12.108578
$ws.Range("O3").Value = 0.9574892080033358
$ws.Range("P3").Value = 0.9574892080033357
$ws.Range("Q3").Value = 34.37578205285556
$ws.Range("R3").Value = 309.3820384757
$ws.Range("S3").Value = 0.08113143414512429
$ws.Range("T3").Value = 0.08113143414512426

# Row 4: M2 -> Cd86/Cd28 -> ECs
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Cd86"
$ws.Range("C4").Value = "Cd28"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 91.99685666666666
$ws.Range("H4").Value = 275.99057
$ws.Range("I4").Value = 0.9152664766694251
$ws.Range("J4").Value = 0.9152664766694251
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1791996666666667
$ws.Range("N4").Value = 0.537599
$ws.Range("O4").Value = 0.04251079199666429
$ws.Range("P4").Value = 0.04251079199666429
$ws.Range("Q4").Value = 16.48580604904778
$ws.Range("R4").Value = 148.37225444143
$ws.Range("S4").Value = 0.03890870281121372
$ws.Range("T4").Value = 0.03890870281121372

# Row 5: M2 -> Cd86/Cd28 -> M2
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Cd86"
$ws.Range("C5").Value = "Cd28"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 91.99685666666666
$ws.Range("H5").Value = 275.99057
$ws.Range("I5").Value = 0.9152664766694251
$ws.Range("J5").Value = 0.9152664766694251
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.036192666666667
$ws.Range("N5").Value = 12.108578
$ws.Range("O5").Value = 0.9574892080033358
$ws.Range("P5").Value = 0.9574892080033357
$ws.Range("Q5").Value = 371.3170382343844
$ws.Range("R5").Value = 3341.85334410946
$ws.Range("S5").Value = 0.8763577738582115
$ws.Range("T5").Value = 0.8763577738582113
